# Insert two new columns (INDIVIDUAL_REF_ID / INDIVIDUAL_REF_DB) right after
# TREATMENT (column R) on every worksheet of the submission template, shifting
# all the subsequent metadata columns two places to the right.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Insert 2 blank columns at S:T - everything from column S onward shifts
    # right by two columns (S,T -> U,V ; ... ; BJ -> BL).
    $ws.Range("S1:T1").EntireColumn.Insert()

    # The newly inserted columns don't inherit the sheet's uniform 20-wide
    # column formatting - set it explicitly to match every other column.
    $ws.Range("S1:T1").ColumnWidth = 19.17

    # New header row labels for the two inserted columns.
    $ws.Cells.Item(1, 19).Value = "INDIVIDUAL_REF_ID"
    $ws.Cells.Item(1, 20).Value = "INDIVIDUAL_REF_DB"

    # The column that used to be IMAGING_DATASET (old AM, now AO after the
    # shift) is renamed to IMAGING_DATASET_ID.
    if ($ws.Cells.Item(1, 41).Value() -eq "IMAGING_DATASET") {
        $ws.Cells.Item(1, 41).Value = "IMAGING_DATASET_ID"
    }

    # The "Examples & Info" sheet also carries description / example /
    # used_for / category / regex rows beneath the header - populate the two
    # newly inserted columns for each of them.
    if ($ws.Cells.Item(2, 19).Value() -ne $null -or $ws.Cells.Item(3,1).Value() -eq "example") {
        # row 2: field descriptions
        $ws.Cells.Item(2, 19).Value = "ID/Label of the individual as referenced to in external database"
        $ws.Cells.Item(2, 20).Value = "Label of the external mouse database e.g. MoVi / CRUK-CI"

        # row 3: example values
        $ws.Cells.Item(3, 19).Value = "[12345]"
        $ws.Cells.Item(3, 20).Value = "['MoVi, CRUK-CI']"

        # row 4: used_for
        $ws.Cells.Item(4, 19).Value = "Odomlab"
        $ws.Cells.Item(4, 20).Value = "Odomlab"

        # row 5: category
        $ws.Cells.Item(5, 19).Value = "sample"
        $ws.Cells.Item(5, 20).Value = "sample"

        # row 6 (regex) is left blank for these two columns, matching the
        # rest of the un-validated fields.
    }
}
